$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert a new column A (shifts old A..D into B..E) ---
$ws.Columns.Item(1).Insert()

# --- Step 2: insert two new rows below row 2 for SamplesTab/FilesTab ---
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

# --- Step 3: header row (row 1) ---
$ws.Range("A1").Value = "TabName"
# B1/C1/E1 already correct (query/StatQuery/WebExcel) from the column shift
# D1 target content (per source edit) equals the StatQuery query text, not "dbExcel"
$statQueryText = @'
MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Bernese Mountain Dog']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study
'@
$ws.Range("D1").Value = $statQueryText

# --- Step 4: row 2 (CasesTab) ---
$ws.Range("A2").Value = "CasesTab"
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Bernese Mountain Dog']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
        coalesce(co.cohort_description, '') AS `Cohort`

'@
$ws.Range("B2").Value = $casesQuery
# C2/D2/E2 already correct (StatQuery query / Neo4jData / WebData) from the column shift

# --- Step 5: row 3 (SamplesTab) ---
$ws.Range("A3").Value = "SamplesTab"
$samplesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE demo.breed IN  ['Bernese Mountain Dog']
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
'@
$ws.Range("B3").Value = $samplesQuery
$ws.Range("C3").Value = $statQueryText
$ws.Range("D3").Value = "TC07_Canine_Filter_Breed-BrnMtnDog_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC07_Canine_Filter_Breed-BrnMtnDog_WebData.xlsx"

# --- Step 6: row 4 (FilesTab) ---
$ws.Range("A4").Value = "FilesTab"
$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN  ['Bernese Mountain Dog']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_type, '') AS `File Type`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@
$ws.Range("B4").Value = $filesQuery
$ws.Range("C4").Value = $statQueryText
$ws.Range("D4").Value = "TC07_Canine_Filter_Breed-BrnMtnDog_Neo4jData.xlsx"
$ws.Range("E4").Value = "TC07_Canine_Filter_Breed-BrnMtnDog_WebData.xlsx"

# --- Step 7: wrap text styling on the query columns (B/C for rows 2-4) ---
$ws.Range("B2:C4").WrapText = $true

# --- Step 8: row heights ---
$ws.Rows.Item(2).RowHeight = 275.5
$ws.Rows.Item(3).RowHeight = 232
$ws.Rows.Item(4).RowHeight = 246.5

# --- Step 9: column widths ---
$ws.Columns.Item(1).ColumnWidth = 10

# --- Step 10: sheet view / zoom ---
$ws.Application.ActiveWindow.Zoom = 70

